# Updates the cryptos list (price + 1h volume% columns, plus the
# Bittensor / PancakeSwap row swap) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.348.74'
$ws.Range("E2").Value = '  +2.58%  '

# Row 3
$ws.Range("D3").Value = '2.656.16'
$ws.Range("E3").Value = '  +1.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''605.83'
$ws.Range("E5").Value = '  +2.02%  '

# Row 6
$ws.Range("D6").Value = '''157.32'
$ws.Range("E6").Value = '  +4.34%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").Value = '''0.588'
$ws.Range("E8").Value = '  +0.47%  '

# Row 9
$ws.Range("E9").Value = '  +8.65%  '

# Row 10
$ws.Range("D10").Value = '''0.413'
$ws.Range("E10").Value = '  +4.90%  '

# Row 11
$ws.Range("D11").Value = '''5.82'
$ws.Range("E11").Value = '  +0.51%  '

# Row 12
$ws.Range("E12").Value = '  +1.61%  '

# Row 13
$ws.Range("D13").Value = '''29.53'
$ws.Range("E13").Value = '  +5.97%  '

# Row 14
$ws.Range("E14").Value = '  +15.47%  '

# Row 15
$ws.Range("D15").Value = '3.133.85'
$ws.Range("E15").Value = '  +1.51%  '

# Row 16
$ws.Range("D16").Value = '65.163.98'
$ws.Range("E16").Value = '  +2.55%  '

# Row 17
$ws.Range("D17").Value = '2.670.89'
$ws.Range("E17").Value = '  +1.77%  '

# Row 18
$ws.Range("D18").Value = '''12.80'
$ws.Range("E18").Value = '  +4.65%  '

# Row 19
$ws.Range("D19").Value = '''4.91'
$ws.Range("E19").Value = '  +2.67%  '

# Row 20
$ws.Range("D20").Value = '''359.72'
$ws.Range("E20").Value = '  +3.19%  '

# Row 21
$ws.Range("D21").Value = '''7.36'
$ws.Range("E21").Value = '  +5.48%  '

# Row 22
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").Value = '''69.17'
$ws.Range("E23").Value = '  +2.85%  '

# Row 24
$ws.Range("E24").Value = '  +2.15%  '

# Row 25
$ws.Range("D25").Value = '''9.53'
$ws.Range("E25").Value = '  +3.27%  '

# Row 26
$ws.Range("D26").Value = '''0.0000103'
$ws.Range("E26").Value = '  +15.60%  '

# Row 27
$ws.Range("E27").Value = '  -1.15%  '

# Row 28
$ws.Range("D28").Value = '''8.26'
$ws.Range("E28").Value = '  -1.85%  '

# Row 29
$ws.Range("E29").Value = '  +1.77%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.24'
$ws.Range("E30").Value = '  +8.56%  '

# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '''547.35'
$ws.Range("E31").Value = '  +0.40%  '

# Row 32
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("D33").Value = '''1.84'
$ws.Range("E33").Value = '  +2.67%  '

# Row 34
$ws.Range("E34").Value = '  +3.62%  '

# Row 35
$ws.Range("D35").Value = '''6.42'
$ws.Range("E35").Value = '  +4.52%  '

# Row 36
$ws.Range("D36").Value = '''0.433'
$ws.Range("E36").Value = '  +3.60%  '

# Row 37
$ws.Range("D37").Value = '''20.64'
$ws.Range("E37").Value = '  +4.64%  '

# Row 38
$ws.Range("D38").Value = '''163.16'
$ws.Range("E38").Value = '  -0.80%  '

# Row 39
$ws.Range("E39").Value = '  +1.74%  '

# Row 40
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.13%  '

# Row 41
$ws.Range("E41").Value = '  +0.05%  '

# Row 42
$ws.Range("D42").Value = '''42.61'
$ws.Range("E42").Value = '  +7.00%  '

# Row 43
$ws.Range("D43").Value = '''167.07'
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("D44").Value = '''4.20'
$ws.Range("E44").Value = '  +2.89%  '

# Row 45
$ws.Range("D45").Value = '''0.0619'
$ws.Range("E45").Value = '  +6.01%  '

# Row 46
$ws.Range("D46").Value = '''2.31'
$ws.Range("E46").Value = '  +7.13%  '

# Row 47
$ws.Range("D47").Value = '''23.15'
$ws.Range("E47").Value = '  -1.23%  '

# Row 48
$ws.Range("E48").Value = '  +3.29%  '

# Row 49
$ws.Range("D49").Value = '''0.0264'

# Row 50
$ws.Range("D50").Value = '''0.0985'
$ws.Range("E50").Value = '  +2.05%  '

# Row 51
$ws.Range("D51").Value = '''19.72'
$ws.Range("E51").Value = '  +2.36%  '
